function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws "D2" "317.81"
Set-TextValue $ws "E2" "3.78%"
Set-TextValue $ws "G2" "18"

# Row 3
Set-TextValue $ws "D3" "39.79"
Set-TextValue $ws "E3" "2.05%"
Set-TextValue $ws "G3" "18"

# Row 4
Set-TextValue $ws "D4" "5.146"
Set-TextValue $ws "E4" "0.91%"
Set-TextValue $ws "G4" "18"

# Row 5
Set-TextValue $ws "D5" "0.08192"
Set-TextValue $ws "E5" "1.71%"
Set-TextValue $ws "G5" "18"

# Row 6
Set-TextValue $ws "D6" "2.081"
Set-TextValue $ws "E6" "6.82%"
Set-TextValue $ws "G6" "18"

# Row 7
Set-TextValue $ws "D7" "8.316"
Set-TextValue $ws "E7" "3.93%"
Set-TextValue $ws "G7" "18"

# Row 8
Set-TextValue $ws "D8" "4.299"
Set-TextValue $ws "E8" "2.50%"
Set-TextValue $ws "G8" "18"

# Row 9
Set-TextValue $ws "D9" "0.9353"
Set-TextValue $ws "E9" "0.43%"
Set-TextValue $ws "G9" "18"

# Row 10
Set-TextValue $ws "D10" "0.1298"
Set-TextValue $ws "E10" "-10.72%"
Set-TextValue $ws "G10" "18"

# Row 11
Set-TextValue $ws "D11" "0.1981"
Set-TextValue $ws "E11" "2.83%"
Set-TextValue $ws "G11" "18"

# Row 12
Set-TextValue $ws "D12" "0.09080"
Set-TextValue $ws "E12" "0.52%"
Set-TextValue $ws "G12" "18"

# Row 13
Set-TextValue $ws "D13" "0.03489"
Set-TextValue $ws "E13" "-0.47%"
Set-TextValue $ws "G13" "18"

# Row 14
Set-TextValue $ws "D14" "0.09826"
Set-TextValue $ws "E14" "0.41%"
Set-TextValue $ws "G14" "18"

# Row 15
Set-TextValue $ws "D15" "0.001409"
Set-TextValue $ws "E15" "0.70%"
Set-TextValue $ws "G15" "18"

# Row 16
Set-TextValue $ws "D16" "0.006379"
Set-TextValue $ws "E16" "6.56%"
Set-TextValue $ws "G16" "18"

# Row 17
Set-TextValue $ws "D17" "3.678"
Set-TextValue $ws "E17" "-2.87%"
Set-TextValue $ws "G17" "18"

# Row 18
Set-TextValue $ws "D18" "3.207"
Set-TextValue $ws "E18" "-6.93%"
Set-TextValue $ws "G18" "18"

# Row 19
Set-TextValue $ws "G19" "18"

# Row 20
Set-TextValue $ws "D20" "0.1293"
Set-TextValue $ws "E20" "-0.75%"
Set-TextValue $ws "G20" "18"

# Row 21
Set-TextValue $ws "E21" "2.39%"
Set-TextValue $ws "G21" "18"

# Row 22
Set-TextValue $ws "D22" "0.2452"
Set-TextValue $ws "E22" "1.57%"
Set-TextValue $ws "G22" "18"

# Row 23
Set-TextValue $ws "D23" "0.04334"
Set-TextValue $ws "E23" "-1.04%"
Set-TextValue $ws "G23" "18"

# Row 24
Set-TextValue $ws "D24" "0.001227"
Set-TextValue $ws "E24" "-0.90%"
Set-TextValue $ws "G24" "18"

# Row 25
Set-TextValue $ws "D25" "0.004751"
Set-TextValue $ws "E25" "11.08%"
Set-TextValue $ws "G25" "18"

# Row 26
Set-TextValue $ws "G26" "18"

# Row 27
Set-TextValue $ws "D27" "0.0004002"
Set-TextValue $ws "E27" "-10.01%"
Set-TextValue $ws "G27" "18"

# Row 28
Set-TextValue $ws "G28" "18"

# Row 29
Set-TextValue $ws "G29" "18"

# Row 30
Set-TextValue $ws "G30" "18"

# Row 31
Set-TextValue $ws "G31" "18"

# Row 32
Set-TextValue $ws "G32" "18"

# Row 33
Set-TextValue $ws "G33" "18"

# Row 34
Set-TextValue $ws "G34" "18"

# Row 35
Set-TextValue $ws "G35" "18"

# Row 36
Set-TextValue $ws "G36" "18"

# Row 37
Set-TextValue $ws "G37" "18"

# Row 38
Set-TextValue $ws "G38" "18"

# Row 39
Set-TextValue $ws "D39" "0.02216"
Set-TextValue $ws "E39" "9.22%"
Set-TextValue $ws "G39" "18"

# Row 40
Set-TextValue $ws "D40" "0.05220"
Set-TextValue $ws "E40" "3.67%"
Set-TextValue $ws "G40" "18"

# Row 41
Set-TextValue $ws "D41" "0.007634"
Set-TextValue $ws "E41" "2.43%"
Set-TextValue $ws "G41" "18"

# Row 42
Set-TextValue $ws "D42" "0.009703"
Set-TextValue $ws "E42" "-4.93%"
Set-TextValue $ws "G42" "18"

# Row 43
Set-TextValue $ws "D43" "0.1379"
Set-TextValue $ws "E43" "2.40%"
Set-TextValue $ws "G43" "18"

# Row 44
Set-TextValue $ws "D44" "0.002130"
Set-TextValue $ws "E44" "0.36%"
Set-TextValue $ws "G44" "18"

# Row 45
Set-TextValue $ws "D45" "0.009188"
Set-TextValue $ws "E45" "0.82%"
Set-TextValue $ws "G45" "18"

# Row 46
Set-TextValue $ws "D46" "0.00006547"
Set-TextValue $ws "E46" "5.75%"
Set-TextValue $ws "G46" "18"

# Row 47
Set-TextValue $ws "G47" "18"

# Row 48
Set-TextValue $ws "B48" "BOLO"
Set-TextValue $ws "C48" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws "D48" "0.002777"
Set-TextValue $ws "E48" "-6.37%"
Set-TextValue $ws "G48" "18"

# Row 49
Set-TextValue $ws "B49" "CoinbaseStockToken"
Set-TextValue $ws "C49" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws "D49" "0.001200"
Set-TextValue $ws "E49" "-25.07%"
Set-TextValue $ws "G49" "18"

# Row 50
Set-TextValue $ws "G50" "18"

# Row 51
Set-TextValue $ws "G51" "18"
